# Append the new daily row (row 48) to Sheet1, extending the data range
# from A1:D47 to A1:D48.
#
# A48 must stay a literal text string ("2025/10/02"), not an Excel date
# serial number, so we flip the cell to Text format before assigning the
# value and then clear the formatting again afterwards so the cell is left
# with the workbook's default (unstyled) look, matching the rest of the
# date column (A2:A47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48

$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/02"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = "木"
$ws.Cells.Item($row, 3).Value = 7
$ws.Cells.Item($row, 4).Value = 24
